$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("A2").Value = "30 Oct 2025, 12:05 PM"

$ws = $wb.Worksheets.Item("Top Gainers")
$ws.Range("C2").Value = 11.4337
$ws.Range("D2").Value = 10.3656
$ws.Range("E2").Value = -19.2572
$ws.Range("C3").Value = 11.1629
$ws.Range("D3").Value = 19.9682
$ws.Range("E3").Value = 27.125
$ws.Range("C4").Value = 11.0085
$ws.Range("D4").Value = 18.8342
$ws.Range("E4").Value = 32.3722
$ws.Range("C5").Value = 9.975099999999999
$ws.Range("D5").Value = 10.1806
$ws.Range("E5").Value = 24.1793
$ws.Range("C6").Value = 9.5405
$ws.Range("D6").Value = 6.6116
$ws.Range("E6").Value = -8.1274
$ws.Range("B7").Value = "SALASAR"
$ws.Range("C7").Value = 8.933999999999999
$ws.Range("D7").Value = 14.1489
$ws.Range("E7").Value = 20.9696
$ws.Range("B8").Value = "UNIPARTS"
$ws.Range("C8").Value = 8.4069
$ws.Range("D8").Value = 10.6986
$ws.Range("E8").Value = 26.9974
$ws.Range("C9").Value = 7.4781
$ws.Range("D9").Value = 12.6756
$ws.Range("E9").Value = 14.839
$ws.Range("B10").Value = "SHREEJISPG"
$ws.Range("C10").Value = 6.3563
$ws.Range("D10").Value = 10.6404
$ws.Range("E10").Value = 11.7378
$ws.Range("B11").Value = "MARINE"
$ws.Range("C11").Value = 6.3006
$ws.Range("D11").Value = 2.831
$ws.Range("E11").Value = 15.3882
$ws.Range("B12").Value = "INDIACEM"
$ws.Range("C12").Value = 6.0687
$ws.Range("D12").Value = 6.3541
$ws.Range("E12").Value = 8.099
$ws.Range("B13").Value = "PDSL"
$ws.Range("C13").Value = 5.7238
$ws.Range("D13").Value = 9.0733
$ws.Range("E13").Value = 15.2444
$ws.Range("C14").Value = 5.6937
$ws.Range("D14").Value = 19.0875
$ws.Range("E14").Value = 19.1763
$ws.Range("C16").Value = 5.3892
$ws.Range("D16").Value = 6.0184
$ws.Range("E16").Value = 3.7124
$ws.Range("C17").Value = 5.1737
$ws.Range("D17").Value = 5.7265
$ws.Range("E17").Value = -0.6175
$ws.Range("C18").Value = 5.0496
$ws.Range("D18").Value = 15.2502
$ws.Range("E18").Value = 26.1164
$ws.Range("B21").Value = "BHEL"
$ws.Range("C21").Value = 4.9554
$ws.Range("D21").Value = 11.5225
$ws.Range("E21").Value = 8.0101
$ws.Range("B22").Value = "VIMTALABS"
$ws.Range("C22").Value = 4.8949
$ws.Range("D22").Value = 5.0744
$ws.Range("E22").Value = -0.0566
$ws.Range("B23").Value = "CREDITACC"
$ws.Range("C23").Value = 4.7304
$ws.Range("D23").Value = 3.3462
$ws.Range("E23").Value = 8.6388
$ws.Range("B24").Value = "VSTIND"
$ws.Range("C24").Value = 4.7186
$ws.Range("D24").Value = 5.2274
$ws.Range("E24").Value = 4.7591
$ws.Range("B25").Value = "POLICYBZR"
$ws.Range("C25").Value = 4.6768
$ws.Range("D25").Value = 7.0179
$ws.Range("E25").Value = 5.9929
$ws.Range("B26").Value = "IIFL"
$ws.Range("C26").Value = 4.6176
$ws.Range("D26").Value = 11.7089
$ws.Range("E26").Value = 21.0777
$ws.Range("B27").Value = "NETWEB"
$ws.Range("C27").Value = 4.3754
$ws.Range("D27").Value = 10.1907
$ws.Range("E27").Value = 12.2626
$ws.Range("B29").Value = "BEML"
$ws.Range("C29").Value = 4.2388
$ws.Range("D29").Value = 1.2888
$ws.Range("E29").Value = 7.7028
$ws.Range("B30").Value = "SHANTIGOLD"
$ws.Range("C30").Value = 4.2168
$ws.Range("D30").Value = 11.5519
$ws.Range("E30").Value = 4.1096
$ws.Range("B31").Value = "EUROPRATIK"
$ws.Range("C31").Value = 4.1753
$ws.Range("D31").Value = 10.0401
$ws.Range("E31").Value = 26.9727
$ws.Range("B32").Value = "DEEDEV"
$ws.Range("C32").Value = 4.1101
$ws.Range("D32").Value = -2.8161
$ws.Range("E32").Value = -3.6176
$ws.Range("B34").Value = "ALICON"
$ws.Range("C34").Value = 3.7111
$ws.Range("D34").Value = 9.8592
$ws.Range("E34").Value = 15.3162
$ws.Range("B35").Value = "JKTYRE"
$ws.Range("C35").Value = 3.7015
$ws.Range("D35").Value = 6.7727
$ws.Range("E35").Value = 22.9834
$ws.Range("B36").Value = "BLS"
$ws.Range("C36").Value = 3.6749
$ws.Range("D36").Value = 0.6127
$ws.Range("E36").Value = -0.6516
$ws.Range("B37").Value = "SHRINGARMS"
$ws.Range("C37").Value = 3.6388
$ws.Range("D37").Value = 4.8508
$ws.Range("E37").Value = 25.1046
$ws.Range("B38").Value = "SKYGOLD"
$ws.Range("C38").Value = 3.6332
$ws.Range("D38").Value = -0.9163
$ws.Range("E38").Value = 37.6068
$ws.Range("B39").Value = "SUNDROP"
$ws.Range("C39").Value = 3.6303
$ws.Range("D39").Value = 3.5224
$ws.Range("E39").Value = 1.6104
$ws.Range("B40").Value = "HIRECT"
$ws.Range("C40").Value = 3.5571
$ws.Range("D40").Value = 10.9669
$ws.Range("E40").Value = 9.2629
$ws.Range("B41").Value = "OIL"
$ws.Range("C41").Value = 3.5315
$ws.Range("D41").Value = 3.7783
$ws.Range("E41").Value = 5.2078
$ws.Range("B42").Value = "MTARTECH"
$ws.Range("C42").Value = 3.4424
$ws.Range("D42").Value = 7.5611
$ws.Range("E42").Value = 31.4681
$ws.Range("B43").Value = "RELTD"
$ws.Range("C43").Value = 3.4353
$ws.Range("D43").Value = 10.2373
$ws.Range("E43").Value = -1.3111
$ws.Range("B44").Value = "RAMASTEEL"
$ws.Range("C44").Value = 3.4137
$ws.Range("D44").Value = 3.3099
$ws.Range("E44").Value = 4.888
$ws.Range("B45").Value = "PSPPROJECT"
$ws.Range("C45").Value = 3.3259
$ws.Range("D45").Value = 17.4617
$ws.Range("E45").Value = 23.9197
$ws.Range("C46").Value = 3.2412
$ws.Range("D46").Value = 4.7863
$ws.Range("B48").Value = "RSYSTEMS"
$ws.Range("C48").Value = 3.0888
$ws.Range("D48").Value = 4.246
$ws.Range("E48").Value = 6.5868
$ws.Range("B49").Value = "SAPPHIRE"
$ws.Range("C49").Value = 3.088
$ws.Range("D49").Value = 4.9057
$ws.Range("E49").Value = 2.2635
$ws.Range("B51").Value = "CENTUM"
$ws.Range("C51").Value = 3.0057
$ws.Range("D51").Value = 3.8656
$ws.Range("E51").Value = -1.6411
$ws.Range("B52").Value = "DBCORP"
$ws.Range("C52").Value = 2.8808
$ws.Range("D52").Value = 5.5234
$ws.Range("E52").Value = 1.467
$ws.Range("B53").Value = "POWERINDIA"
$ws.Range("C53").Value = 2.8601
$ws.Range("D53").Value = 7.2464
$ws.Range("E53").Value = -0.1055
$ws.Range("B54").Value = "AHLUCONT"
$ws.Range("C54").Value = 2.8559
$ws.Range("D54").Value = 1.6656
$ws.Range("E54").Value = -5.5466
$ws.Range("B56").Value = "CARYSIL"
$ws.Range("C56").Value = 2.8151
$ws.Range("D56").Value = 2.2987
$ws.Range("E56").Value = 11.1995
$ws.Range("B58").Value = "NEULANDLAB"
$ws.Range("C58").Value = 2.7652
$ws.Range("D58").Value = -1.5912
$ws.Range("E58").Value = 8.401199999999999
$ws.Range("B59").Value = "OBEROIRLTY"
$ws.Range("C59").Value = 2.7256
$ws.Range("D59").Value = 3.5414
$ws.Range("E59").Value = 11.2509
$ws.Range("C60").Value = 2.6992
$ws.Range("D60").Value = 5.945
$ws.Range("E60").Value = 10.5062
$ws.Range("B61").Value = "GREAVESCOT"
$ws.Range("C61").Value = 2.6823
$ws.Range("D61").Value = 15.2532
$ws.Range("E61").Value = 11.494
$ws.Range("B62").Value = "MAMATA"
$ws.Range("C62").Value = 2.6048
$ws.Range("D62").Value = 2.0155
$ws.Range("E62").Value = 1.2177
$ws.Range("B63").Value = "CIFL"
$ws.Range("C63").Value = 2.588
$ws.Range("D63").Value = 2.1424
$ws.Range("E63").Value = 2.0833
$ws.Range("B64").Value = "DIVISLAB"
$ws.Range("C64").Value = 2.5799
$ws.Range("D64").Value = 1.3196
$ws.Range("E64").Value = 17.4093
$ws.Range("B65").Value = "DBL"
$ws.Range("C65").Value = 2.5484
$ws.Range("D65").Value = 3.6327
$ws.Range("E65").Value = 4.707
$ws.Range("B67").Value = "BLISSGVS"
$ws.Range("C67").Value = 2.5189
$ws.Range("D67").Value = 1.8757
$ws.Range("E67").Value = 2.2432
$ws.Range("B68").Value = "BLUEDART"
$ws.Range("C68").Value = 2.5183
$ws.Range("D68").Value = 21.0583
$ws.Range("E68").Value = 18.1396
$ws.Range("B69").Value = "INDOSTAR"
$ws.Range("C69").Value = 2.503
$ws.Range("D69").Value = 7.3646
$ws.Range("E69").Value = 1.9866
$ws.Range("B70").Value = "IVALUE"
$ws.Range("C70").Value = 2.4972
$ws.Range("D70").Value = 5.7286
$ws.Range("E70").Value = -1.4337
$ws.Range("B71").Value = "CAMS"
$ws.Range("C71").Value = 2.4199
$ws.Range("D71").Value = 1.768
$ws.Range("E71").Value = 5.0604
$ws.Range("B72").Value = "GPPL"
$ws.Range("C72").Value = 2.4035
$ws.Range("D72").Value = 5.8926
$ws.Range("E72").Value = 7.5745
$ws.Range("B73").Value = "JSFB"
$ws.Range("C73").Value = 2.3807
$ws.Range("D73").Value = 2.2429
$ws.Range("E73").Value = -4.6937
$ws.Range("B74").Value = "CEATLTD"
$ws.Range("C74").Value = 2.3463
$ws.Range("D74").Value = -0.7704
$ws.Range("E74").Value = 21.4459
$ws.Range("B75").Value = "REDTAPE"
$ws.Range("C75").Value = 2.3237
$ws.Range("D75").Value = 2.2244
$ws.Range("E75").Value = -4.6176
$ws.Range("B76").Value = "EPACK"
$ws.Range("C76").Value = 2.3007
$ws.Range("D76").Value = 1.6453
$ws.Range("E76").Value = 0.5763

$ws = $wb.Worksheets.Item("Top Losers")
$ws.Range("C2").Value = -18.5527
$ws.Range("D2").Value = -17.2921
$ws.Range("E2").Value = -0.264
$ws.Range("C3").Value = -9.334
$ws.Range("D3").Value = -5.8888
$ws.Range("E3").Value = 6.376
$ws.Range("C4").Value = -5.8761
$ws.Range("D4").Value = -8.42
$ws.Range("E4").Value = 8.364100000000001
$ws.Range("C10").Value = -4.6
$ws.Range("D10").Value = -2.9352
$ws.Range("E10").Value = 20.2036
$ws.Range("C11").Value = -4.5967
$ws.Range("D11").Value = -7.0271
$ws.Range("E11").Value = -2.4761
$ws.Range("C12").Value = -4.3488
$ws.Range("D12").Value = -3.1018
$ws.Range("E12").Value = -3.901
$ws.Range("C13").Value = -4.275
$ws.Range("D13").Value = -3.4909
$ws.Range("E13").Value = 21.3821
$ws.Range("C15").Value = -3.596
$ws.Range("D15").Value = -6.3391
$ws.Range("E15").Value = 3.51
$ws.Range("B16").Value = "TVSHLTD"
$ws.Range("C16").Value = -3.4813
$ws.Range("D16").Value = -2.2385
$ws.Range("E16").Value = 16.0266
$ws.Range("C17").Value = -3.4777
$ws.Range("D17").Value = -3.2791
$ws.Range("E17").Value = 5.1508
$ws.Range("B18").Value = "LICHSGFIN"
$ws.Range("C18").Value = -3.4712
$ws.Range("D18").Value = -1.1475
$ws.Range("E18").Value = 1.3894
$ws.Range("C19").Value = -3.3743
$ws.Range("D19").Value = 4.1615
$ws.Range("E19").Value = -14.1645
$ws.Range("C20").Value = -3.3582
$ws.Range("D20").Value = 4.9204
$ws.Range("E20").Value = 1.0039
$ws.Range("B21").Value = "QUESS"
$ws.Range("C21").Value = -3.3028
$ws.Range("D21").Value = 3.8018
$ws.Range("E21").Value = -4.5562
$ws.Range("B22").Value = "UBL"
$ws.Range("C22").Value = -3.2766
$ws.Range("D22").Value = -2.6844
$ws.Range("E22").Value = -1.155
$ws.Range("B23").Value = "FILATEX"
$ws.Range("C23").Value = -3.2764
$ws.Range("D23").Value = 6.661
$ws.Range("E23").Value = 21.8744
$ws.Range("B24").Value = "BHARTIHEXA"
$ws.Range("C24").Value = -3.2649
$ws.Range("D24").Value = 3.5914
$ws.Range("E24").Value = 11.5677
$ws.Range("B26").Value = "APOLLOPIPE"
$ws.Range("C26").Value = -3.2591
$ws.Range("D26").Value = -4.9461
$ws.Range("E26").Value = -9.970000000000001
$ws.Range("B27").Value = "CGCL"
$ws.Range("C27").Value = -3.2236
$ws.Range("D27").Value = -1.2783
$ws.Range("E27").Value = 9.2667
$ws.Range("B28").Value = "RAMCOSYS"
$ws.Range("C28").Value = -3.1844
$ws.Range("D28").Value = 6.6448
$ws.Range("E28").Value = 25.2683
$ws.Range("B29").Value = "INDUSTOWER"
$ws.Range("C29").Value = -3.1623
$ws.Range("D29").Value = 2.0606
$ws.Range("E29").Value = 7.6116
$ws.Range("B31").Value = "EPACKPEB"
$ws.Range("C31").Value = -2.9953
$ws.Range("D31").Value = -2.9303
$ws.Range("E31").Value = "N/A"
$ws.Range("B32").Value = "SOUTHBANK"
$ws.Range("C32").Value = -2.966
$ws.Range("D32").Value = -0.4198
$ws.Range("E32").Value = 31.0881
$ws.Range("B33").Value = "BCG"
$ws.Range("C33").Value = -2.9161
$ws.Range("D33").Value = 2.0942
$ws.Range("E33").Value = -1.7279
$ws.Range("B34").Value = "UTIAMC"
$ws.Range("C34").Value = -2.8039
$ws.Range("D34").Value = -7.4088
$ws.Range("E34").Value = -4.8386
$ws.Range("B35").Value = "BHARATWIRE"
$ws.Range("C35").Value = -2.7123
$ws.Range("D35").Value = 19.502
$ws.Range("E35").Value = 20.5374
$ws.Range("C36").Value = -2.7092
$ws.Range("D36").Value = -3.1466
$ws.Range("E36").Value = -7.0129
$ws.Range("B37").Value = "POCL"
$ws.Range("C37").Value = -2.6978
$ws.Range("D37").Value = 2.4547
$ws.Range("E37").Value = 22.9331
$ws.Range("B38").Value = "SURAJEST"
$ws.Range("C38").Value = -2.6742
$ws.Range("D38").Value = 6.2979
$ws.Range("E38").Value = 4.2988
$ws.Range("B39").Value = "DREDGECORP"
$ws.Range("C39").Value = -2.6671
$ws.Range("D39").Value = 18.515
$ws.Range("E39").Value = 19.2897
$ws.Range("B41").Value = "TVSELECT"
$ws.Range("C41").Value = -2.6305
$ws.Range("D41").Value = -3.5787
$ws.Range("E41").Value = -5.5485
$ws.Range("B42").Value = "IDEAFORGE"
$ws.Range("C42").Value = -2.6234
$ws.Range("D42").Value = -1.8173
$ws.Range("E42").Value = -3.5697
$ws.Range("B43").Value = "CAMLINFINE"
$ws.Range("C43").Value = -2.6128
$ws.Range("D43").Value = 0.1686
$ws.Range("E43").Value = 0.4303
$ws.Range("B44").Value = "PRECWIRE"
$ws.Range("C44").Value = -2.6124
$ws.Range("D44").Value = 9.4214
$ws.Range("E44").Value = 20.0654
$ws.Range("B45").Value = "TTKPRESTIG"
$ws.Range("C45").Value = -2.5915
$ws.Range("D45").Value = 5.2024
$ws.Range("E45").Value = 6.8089
$ws.Range("B46").Value = "LXCHEM"
$ws.Range("C46").Value = -2.569
$ws.Range("D46").Value = -3.0082
$ws.Range("E46").Value = -4.1818
$ws.Range("B47").Value = "SANDHAR"
$ws.Range("C47").Value = -2.5499
$ws.Range("D47").Value = 1.1702
$ws.Range("E47").Value = 18.5568
$ws.Range("B48").Value = "PDMJEPAPER"
$ws.Range("C48").Value = -2.5337
$ws.Range("D48").Value = -2.451
$ws.Range("E48").Value = -4.2298
$ws.Range("B49").Value = "APARINDS"
$ws.Range("C49").Value = -2.4836
$ws.Range("D49").Value = 5.6507
$ws.Range("E49").Value = 12.7168
$ws.Range("B50").Value = "MSPL"
$ws.Range("C50").Value = -2.449
$ws.Range("D50").Value = -1.2979
$ws.Range("E50").Value = -8.328799999999999
$ws.Range("B51").Value = "VGUARD"
$ws.Range("C51").Value = -2.418
$ws.Range("D51").Value = 0.2955
$ws.Range("E51").Value = -0.5859
$ws.Range("B52").Value = "IDBI"
$ws.Range("C52").Value = -2.3943
$ws.Range("D52").Value = 5.898
$ws.Range("E52").Value = 8.8055
$ws.Range("B53").Value = "GRWRHITECH"
$ws.Range("C53").Value = -2.337
$ws.Range("D53").Value = -5.9316
$ws.Range("E53").Value = 18.8333
$ws.Range("B54").Value = "SINDHUTRAD"
$ws.Range("C54").Value = -2.3237
$ws.Range("D54").Value = -1.2366
$ws.Range("E54").Value = -15.1858
$ws.Range("B55").Value = "HONASA"
$ws.Range("C55").Value = -2.3015
$ws.Range("D55").Value = -0.9048
$ws.Range("E55").Value = -2.5276
$ws.Range("B56").Value = "AEROFLEX"
$ws.Range("C56").Value = -2.2725
$ws.Range("D56").Value = 5.0601
$ws.Range("E56").Value = 3.9084
$ws.Range("B57").Value = "GABRIEL"
$ws.Range("C57").Value = -2.2415
$ws.Range("D57").Value = 1.8917
$ws.Range("E57").Value = 6.687
$ws.Range("B58").Value = "JINDALPHOT"
$ws.Range("C58").Value = -2.2381
$ws.Range("D58").Value = -2.4543
$ws.Range("E58").Value = 20.2307
$ws.Range("B59").Value = "DCMSRIND"
$ws.Range("C59").Value = -2.212
$ws.Range("D59").Value = -0.947
$ws.Range("E59").Value = 4.9242
$ws.Range("B60").Value = "VEDL"
$ws.Range("C60").Value = -2.1988
$ws.Range("D60").Value = 1.8664
$ws.Range("E60").Value = 8.3834
$ws.Range("B61").Value = "ARIHANTCAP"
$ws.Range("C61").Value = -2.1862
$ws.Range("D61").Value = 4.7844
$ws.Range("E61").Value = -4.0497
$ws.Range("B62").Value = "JTEKTINDIA"
$ws.Range("C62").Value = -2.1693
$ws.Range("D62").Value = 4.1165
$ws.Range("E62").Value = -1.8815
$ws.Range("B63").Value = "HFCL"
$ws.Range("C63").Value = -2.1646
$ws.Range("D63").Value = -2.9071
$ws.Range("E63").Value = 3.5675
$ws.Range("B64").Value = "HCG"
$ws.Range("C64").Value = -2.1448
$ws.Range("D64").Value = 0.0725
$ws.Range("E64").Value = 18.0864
$ws.Range("B65").Value = "DCBBANK"
$ws.Range("C65").Value = -2.1434
$ws.Range("D65").Value = -1.2621
$ws.Range("E65").Value = 22.8341
$ws.Range("B66").Value = "DELHIVERY"
$ws.Range("C66").Value = -2.1147
$ws.Range("D66").Value = 1.6171
$ws.Range("E66").Value = 5.4216
$ws.Range("B67").Value = "NEWGEN"
$ws.Range("C67").Value = -2.1039
$ws.Range("D67").Value = 9.1907
$ws.Range("E67").Value = 9.565
$ws.Range("B68").Value = "HMT"
$ws.Range("C68").Value = -2.0681
$ws.Range("D68").Value = -2.5592
$ws.Range("E68").Value = -5.9579
$ws.Range("B69").Value = "SSWL"
$ws.Range("C69").Value = -2.0443
$ws.Range("D69").Value = 4.1651
$ws.Range("E69").Value = 1.2654
$ws.Range("B70").Value = "VIPIND"
$ws.Range("C70").Value = -2.0391
$ws.Range("D70").Value = -3.7718
$ws.Range("E70").Value = -1.2502
$ws.Range("B71").Value = "HDFCLIFE"
$ws.Range("C71").Value = -2.0097
$ws.Range("D71").Value = 1.5035
$ws.Range("E71").Value = -1.3815
$ws.Range("B73").Value = "PROSTARM"
$ws.Range("C73").Value = -1.9857
$ws.Range("D73").Value = -1.327
$ws.Range("E73").Value = -10.0514
$ws.Range("C74").Value = -1.9609
$ws.Range("D74").Value = -6.8865
$ws.Range("E74").Value = -9.5776
$ws.Range("B75").Value = "VBL"
$ws.Range("C75").Value = -1.9578
$ws.Range("D75").Value = 5.2774
$ws.Range("E75").Value = 9.4771
$ws.Range("B76").Value = "TBOTEK"
$ws.Range("C76").Value = -1.9534
$ws.Range("D76").Value = -5.4568
$ws.Range("E76").Value = -0.9376

$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Range("C2").Value = 110.9985
$ws.Range("C4").Value = 78.6481
$ws.Range("C5").Value = 66.45650000000001
$ws.Range("C7").Value = 58.7957
$ws.Range("C9").Value = 52.4108
$ws.Range("C10").Value = 44.1051
$ws.Range("C14").Value = 39.0274
$ws.Range("C15").Value = 38.0833
$ws.Range("C16").Value = 37.0684
$ws.Range("C17").Value = 36.4266
$ws.Range("C18").Value = 35.2687
$ws.Range("C20").Value = 34.4205
$ws.Range("C22").Value = 33.7184
$ws.Range("C23").Value = 32.0354
$ws.Range("C24").Value = 31.4058
$ws.Range("C25").Value = 30.969
$ws.Range("C26").Value = 29.5241
$ws.Range("C28").Value = 28.9685
$ws.Range("B29").Value = "TARACHAND"
$ws.Range("C29").Value = 28.6395
$ws.Range("B30").Value = "ARFIN"
$ws.Range("C30").Value = 28.5781
$ws.Range("C31").Value = 26.8926
$ws.Range("C33").Value = 26.2464
$ws.Range("C34").Value = 25.779
$ws.Range("C35").Value = 25.7727
$ws.Range("C36").Value = 25.0216
$ws.Range("C37").Value = 24.8392
$ws.Range("B38").Value = "MINDTECK"
$ws.Range("C38").Value = 24.5686
$ws.Range("B39").Value = "MARINE"
$ws.Range("C39").Value = 24.4834
$ws.Range("B40").Value = "UNIPARTS"
$ws.Range("C40").Value = 24.4541
$ws.Range("B41").Value = "RAMCOSYS"
$ws.Range("C41").Value = 24.2494
$ws.Range("B42").Value = "TDPOWERSYS"
$ws.Range("C42").Value = 24.175
$ws.Range("C43").Value = 23.8259
$ws.Range("B44").Value = "CARTRADE"
$ws.Range("C44").Value = 23.5834
$ws.Range("B45").Value = "CPEDU"
$ws.Range("C45").Value = 23.5429
$ws.Range("B48").Value = "TATVA"
$ws.Range("C48").Value = 22.7726
$ws.Range("B49").Value = "INDIANB"
$ws.Range("C49").Value = 22.7603
$ws.Range("C50").Value = 22.4871
$ws.Range("B52").Value = "DCBBANK"
$ws.Range("C52").Value = 22.2842
$ws.Range("B53").Value = "KERNEX"
$ws.Range("C53").Value = 22.2181
$ws.Range("C54").Value = 22.1516
$ws.Range("C55").Value = 21.8536
$ws.Range("C56").Value = 21.1803
$ws.Range("B57").Value = "INDRAMEDCO"
$ws.Range("C57").Value = 20.7707
$ws.Range("B58").Value = "PRIVISCL"
$ws.Range("C58").Value = 20.7471
$ws.Range("B59").Value = "SKMEGGPROD"
$ws.Range("C59").Value = 20.7435
$ws.Range("C61").Value = 20.3146
$ws.Range("B62").Value = "SHRIRAMFIN"
$ws.Range("C62").Value = 20.0615
$ws.Range("B63").Value = "BHARATWIRE"
$ws.Range("C63").Value = 20.0245
$ws.Range("C64").Value = 19.9689
$ws.Range("C65").Value = 19.7286
$ws.Range("C66").Value = 19.5582
$ws.Range("C67").Value = 19.5525
$ws.Range("B68").Value = "MCX"
$ws.Range("C68").Value = 19.5213
$ws.Range("B69").Value = "CEATLTD"
$ws.Range("C69").Value = 19.2432
$ws.Range("B70").Value = "REPRO"
$ws.Range("C70").Value = 19.104
$ws.Range("C71").Value = 18.8991
$ws.Range("C72").Value = 18.6837
$ws.Range("C73").Value = 18.4488
$ws.Range("C74").Value = 18.2308
$ws.Range("B76").Value = "THOMASCOTT"
$ws.Range("C76").Value = 18.0348

$ws = $wb.Worksheets.Item("distance from Dma50")
$ws.Range("C2").Value = 9.6488
$ws.Range("C3").Value = 7.259
$ws.Range("C4").Value = 6.5262
$ws.Range("C5").Value = 5.3331
$ws.Range("C6").Value = 5.2867
$ws.Range("C7").Value = 5.0523
$ws.Range("C8").Value = 4.4102
$ws.Range("C9").Value = 4.3705
$ws.Range("C10").Value = 3.876
$ws.Range("C11").Value = 3.7005
$ws.Range("B12").Value = "CNXMIDCAP"
$ws.Range("C12").Value = 3.3969
$ws.Range("B13").Value = "NIFTYFINSERVICE"
$ws.Range("C13").Value = 3.3891
$ws.Range("C14").Value = 3.0494
$ws.Range("C15").Value = 3.0321
$ws.Range("C16").Value = 2.9454
$ws.Range("C17").Value = 2.8265
$ws.Range("C18").Value = 2.8207
$ws.Range("C19").Value = 2.8028
$ws.Range("C20").Value = 2.4017
$ws.Range("C21").Value = 2.2791
$ws.Range("C22").Value = 1.4034
$ws.Range("B23").Value = "CNXIT"
$ws.Range("C23").Value = 1.2995
$ws.Range("B24").Value = "NIFTYCONSUMPTION"
$ws.Range("C24").Value = 1.298
$ws.Range("C25").Value = 1.0377
$ws.Range("C26").Value = 0.9507
$ws.Range("C27").Value = 0.8832
$ws.Range("C28").Value = 0.5323
$ws.Range("C29").Value = 0.4029
$ws.Range("C30").Value = -2.139
